# "Enrollment 4 and 9" - add Grade4StudentCredentials and
# Grade9StudentCredentials worksheets (after GradeOneStudentCredentials),
# each with a credentials header row plus the newly auto-generated
# accounts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# The previously-active sheet loses its single-cell selection/tabSelected
# state once a later sheet becomes active - it is left selecting its
# header row instead.
# ---------------------------------------------------------------------
$gradeOne = $wb.Worksheets.Item("GradeOneStudentCredentials")
$gradeOne.Range("A1:D1").Select()

# ---------------------------------------------------------------------
# New sheet: Grade4StudentCredentials (appended after the last sheet)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$grade4 = $wb.Worksheets.Add($null, $lastSheet)
$grade4.Name = "Grade4StudentCredentials"

$grade4.Range("A1").Value = "UserId"
$grade4.Range("B1").Value = "Password"
$grade4.Range("C1").Value = "UserName"
$grade4.Range("D1").Value = "Signature"
$grade4.Range("A1:B1").Font.Size = 12

$grade4.Range("A2").Value = "Auto2021_07_16_01_25_46_129"
$grade4.Range("B2").Value = "Password@123"

$grade4.Range("A3").Value = "Auto2021_07_16_01_33_33_468"
$grade4.Range("B3").Value = "Password@123"

$grade4.Columns.Item(1).ColumnWidth = 29.6666666666667

$grade4.Range("A1:D1").Select()

# ---------------------------------------------------------------------
# New sheet: Grade9StudentCredentials (appended after Grade4)
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$grade9 = $wb.Worksheets.Add($null, $lastSheet2)
$grade9.Name = "Grade9StudentCredentials"

$grade9.Range("A1").Value = "UserId"
$grade9.Range("B1").Value = "Password"
$grade9.Range("C1").Value = "UserName"
$grade9.Range("D1").Value = "Signature"
$grade9.Range("A1:B1").Font.Size = 12

$grade9.Range("A2").Value = "Auto20210719181533106"
$grade9.Range("B2").Value = "Password@123"

$grade9.Columns.Item(1).ColumnWidth = 6.25

$grade9.Range("G21").Select()
